# Clear the "Included" values in D6:E9 (but keep existing formatting/style)
# and clear the "Included?" value in G23, then move the active selection to D6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6:E6").ClearContents()
$ws.Range("D7:E7").ClearContents()
$ws.Range("D8:E8").ClearContents()
$ws.Range("D9:E9").ClearContents()
$ws.Range("G23").ClearContents()

$ws.Range("D6").Select()
